$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citta")

# Header for new column F
$ws.Range('F1').Value = 'name'

# Citta names - set in the same order they were authored so that
# new shared-strings are appended in a stable, deterministic order
$ws.Range('F2').Value = 'somanassasahagata diṭṭhigatasampayutta asaṅkhārika citta'
$ws.Range('F4').Value = 'somanassasahagata diṭṭhigatavippayutta asaṅkhārika citta'
$ws.Range('F3').Value = 'somanassasahagata diṭṭhigatasampayutta sasaṅkhārika citta'
$ws.Range('F5').Value = 'somanassasahagata diṭṭhigatavippayutta sasaṅkhārika citta'
$ws.Range('F6').Value = 'upekkhāsahagata diṭṭhigatasampayutta asaṅkhārika citta'
$ws.Range('F7').Value = 'upekkhāsahagata diṭṭhigatasampayutta sasaṅkhārika citta'
$ws.Range('F8').Value = 'upekkhāsahagata diṭṭhigatavippayutta asaṅkhārika citta'
$ws.Range('F9').Value = 'upekkhāsahagata diṭṭhigatavippayutta sasaṅkhārika citta'
$ws.Range('F32').Value = 'somanassasahagata ñāṇasampayutta asaṅkhārika kusala citta'
$ws.Range('F33').Value = 'somanassasahagata ñāṇasampayutta sasaṅkhārika kusala citta'
$ws.Range('F34').Value = 'somanassasahagata ñāṇavippayutta asaṅkhārika kusala citta'
$ws.Range('F35').Value = 'somanassasahagata ñāṇavippayutta sasaṅkhārika kusala citta'
$ws.Range('F36').Value = 'upekkhāsahagata ñāṇasampayutta asaṅkhārika kusala citta'
$ws.Range('F37').Value = 'upekkhāsahagata ñāṇasampayutta sasaṅkhārika kusala citta'
$ws.Range('F38').Value = 'upekkhāsahagata ñāṇavippayutta asaṅkhārika kusala citta'
$ws.Range('F39').Value = 'upekkhāsahagata ñāṇavippayutta sasaṅkhārika kusala citta'
$ws.Range('F40').Value = 'somanassasahagata ñāṇasampayutta asaṅkhārika vipāka citta'
$ws.Range('F41').Value = 'somanassasahagata ñāṇasampayutta sasaṅkhārika vipāka citta'
$ws.Range('F42').Value = 'somanassasahagata ñāṇavippayutta asaṅkhārika vipāka citta'
$ws.Range('F43').Value = 'somanassasahagata ñāṇavippayutta sasaṅkhārika vipāka citta'
$ws.Range('F44').Value = 'upekkhāsahagata ñāṇasampayutta asaṅkhārika vipāka citta'
$ws.Range('F45').Value = 'upekkhāsahagata ñāṇasampayutta sasaṅkhārika vipāka citta'
$ws.Range('F46').Value = 'upekkhāsahagata ñāṇavippayutta asaṅkhārika vipāka citta'
$ws.Range('F47').Value = 'upekkhāsahagata ñāṇavippayutta sasaṅkhārika vipāka citta'
$ws.Range('F48').Value = 'somanassasahagata ñāṇasampayutta asaṅkhārika kiriya citta'
$ws.Range('F49').Value = 'somanassasahagata ñāṇasampayutta sasaṅkhārika kiriya citta'
$ws.Range('F50').Value = 'somanassasahagata ñāṇavippayutta asaṅkhārika kiriya citta'
$ws.Range('F51').Value = 'somanassasahagata ñāṇavippayutta sasaṅkhārika kiriya citta'
$ws.Range('F52').Value = 'upekkhāsahagata ñāṇasampayutta asaṅkhārika kiriya citta'
$ws.Range('F53').Value = 'upekkhāsahagata ñāṇasampayutta sasaṅkhārika kiriya citta'
$ws.Range('F54').Value = 'upekkhāsahagata ñāṇavippayutta asaṅkhārika kiriya citta'
$ws.Range('F55').Value = 'upekkhāsahagata ñāṇavippayutta sasaṅkhārika kiriya citta'
$ws.Range('F10').Value = 'domanassasahagata paṭighasampayutta asaṅkhārika citta'
$ws.Range('F11').Value = 'domanassasahagata paṭighasampayutta sasaṅkhārika citta'
$ws.Range('F13').Value = 'upekkhāsahagata uddhacchasampayutta citta'
$ws.Range('F12').Value = 'upekkhāsahagata vicikicchāsampayutta citta'
$ws.Range('F14').Value = 'upekkhāsahagata akusalavipāka cakkhuviññāṇa'
$ws.Range('F15').Value = 'upekkhāsahagata akusalavipāka sotaviññāṇa'
$ws.Range('F16').Value = 'upekkhāsahagata akusalavipāka ghānaviññāṇa'
$ws.Range('F17').Value = 'upekkhāsahagata akusalavipāka jivhāviññāṇa'
$ws.Range('F18').Value = 'dukkhasahagata akusalavipāka kāyaviññāṇa'
$ws.Range('F19').Value = 'upekkhāsahagata akusalavipāka sampaṭicchana citta'
$ws.Range('F20').Value = 'upekkhāsahagata akusalavipāka santīraṇa citta'
$ws.Range('F21').Value = 'upekkhāsahagata kusalavipāka cakkhuviññāṇa'
$ws.Range('F22').Value = 'upekkhāsahagata kusalavipāka sotaviññāṇa'
$ws.Range('F23').Value = 'upekkhāsahagata kusalavipāka ghānaviññāṇa'
$ws.Range('F24').Value = 'upekkhāsahagata kusalavipāka jivhāviññāṇa'
$ws.Range('F25').Value = 'sukhasahagata kusalavipāka kāyaviññāṇa'
$ws.Range('F26').Value = 'upekkhāsahagata kusalavipāka sampaṭicchana citta'
$ws.Range('F27').Value = 'somanassasahagata kusalavipāka santīraṇa citta'
$ws.Range('F28').Value = 'upekkhāsahagata kusalavipāka santīraṇa citta'
$ws.Range('F29').Value = 'upekkhāsahagata pañcadvārāvajjana citta'
$ws.Range('F30').Value = 'upekkhāsahagata manodvārāvajjanacitta citta'
$ws.Range('F31').Value = 'somanassasahagata hasituppāda citta'
$ws.Range('F83').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna sotāpattimagga citta'
$ws.Range('F84').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna sotāpattimagga citta'
$ws.Range('F85').Value = 'pītisukhekaggatāsahita tatiyajjhāna sotāpattimagga citta'
$ws.Range('F86').Value = 'sukhekaggatāsahita catutthajjhāna sotāpattimagga citta'
$ws.Range('F87').Value = 'upekkhekaggatāsahita pañcamajjhāna sotāpattimagga citta'
$ws.Range('F88').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna sakadāgāmimagga citta'
$ws.Range('F89').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna sakadāgāmimagga citta'
$ws.Range('F90').Value = 'pītisukhekaggatāsahita tatiyajjhāna sakadāgāmimagga citta'
$ws.Range('F91').Value = 'sukhekaggatāsahita catutthajjhāna sakadāgāmimagga citta'
$ws.Range('F92').Value = 'upekkhekaggatāsahita pañcamajjhāna sakadāgāmimagga citta'
$ws.Range('F93').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna anāgāmimagga citta'
$ws.Range('F94').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna anāgāmimagga citta'
$ws.Range('F95').Value = 'pītisukhekaggatāsahita tatiyajjhāna anāgāmimagga citta'
$ws.Range('F96').Value = 'sukhekaggatāsahita catutthajjhāna anāgāmimagga citta'
$ws.Range('F97').Value = 'upekkhekaggatāsahita pañcamajjhāna anāgāmimagga citta'
$ws.Range('F98').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna arahattamagga citta'
$ws.Range('F99').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna arahattamagga citta'
$ws.Range('F100').Value = 'pītisukhekaggatāsahita tatiyajjhāna arahattamagga citta'
$ws.Range('F101').Value = 'sukhekaggatāsahita catutthajjhāna arahattamagga citta'
$ws.Range('F102').Value = 'upekkhekaggatāsahita pañcamajjhāna arahattamagga citta'
$ws.Range('F103').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna sotāpattiphala citta'
$ws.Range('F104').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna sotāpattiphala citta'
$ws.Range('F105').Value = 'pītisukhekaggatāsahita tatiyajjhāna sotāpattiphala citta'
$ws.Range('F106').Value = 'sukhekaggatāsahita catutthajjhāna sotāpattiphala citta'
$ws.Range('F107').Value = 'upekkhekaggatāsahita pañcamajjhāna sotāpattiphala citta'
$ws.Range('F108').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna sakadāgāmiphala citta'
$ws.Range('F109').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna sakadāgāmiphala citta'
$ws.Range('F110').Value = 'pītisukhekaggatāsahita tatiyajjhāna sakadāgāmiphala citta'
$ws.Range('F111').Value = 'sukhekaggatāsahita catutthajjhāna sakadāgāmiphala citta'
$ws.Range('F112').Value = 'upekkhekaggatāsahita pañcamajjhāna sakadāgāmiphala citta'
$ws.Range('F113').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna anāgāmiphala citta'
$ws.Range('F114').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna anāgāmiphala citta'
$ws.Range('F115').Value = 'pītisukhekaggatāsahita tatiyajjhāna anāgāmiphala citta'
$ws.Range('F116').Value = 'sukhekaggatāsahita catutthajjhāna anāgāmiphala citta'
$ws.Range('F117').Value = 'upekkhekaggatāsahita pañcamajjhāna anāgāmiphala citta'
$ws.Range('F118').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhāna arahattaphala citta'
$ws.Range('F119').Value = 'vicārapītisukhekaggatāsahita dutiyajjhāna arahattaphala citta'
$ws.Range('F120').Value = 'pītisukhekaggatāsahita tatiyajjhāna arahattaphala citta'
$ws.Range('F121').Value = 'sukhekaggatāsahita catutthajjhāna arahattaphala citta'
$ws.Range('F122').Value = 'upekkhekaggatāsahita pañcamajjhāna arahattaphala citta'
$ws.Range('F56').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhānakusala citta'
$ws.Range('F57').Value = 'vicārapītisukhekaggatāsahita dutiyajjhānakusala citta'
$ws.Range('F58').Value = 'pītisukhekaggatāsahita tatiyajjhānakusala citta'
$ws.Range('F59').Value = 'sukhekaggatāsahita catutthajjhānakusala citta'
$ws.Range('F60').Value = 'upekkhekaggatāsahita pañcamajjhānakusala citta'
$ws.Range('F61').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhānavipāka citta'
$ws.Range('F62').Value = 'vicārapītisukhekaggatāsahita dutiyajjhānavipāka citta'
$ws.Range('F63').Value = 'pītisukhekaggatāsahita tatiyajjhānavipāka citta'
$ws.Range('F64').Value = 'sukhekaggatāsahita catutthajjhānavipāka citta'
$ws.Range('F65').Value = 'upekkhekaggatāsahita pañcamajjhānavipāka citta'
$ws.Range('F66').Value = 'vitakkavicārapītisukhekaggatāsahita paṭhamajjhānakiriya citta'
$ws.Range('F67').Value = 'vicārapītisukhekaggatāsahita dutiyajjhānakiriya citta'
$ws.Range('F68').Value = 'pītisukhekaggatāsahita tatiyajjhānakiriya citta'
$ws.Range('F69').Value = 'sukhekaggatāsahita catutthajjhānakiriya citta'
$ws.Range('F70').Value = 'upekkhekaggatāsahita pañcamajjhānakiriya citta'
$ws.Range('F71').Value = 'ākāsānañcāyatana kusala citta'
$ws.Range('F72').Value = 'viññāṇañcāyatana kusala citta'
$ws.Range('F73').Value = 'ākiñcaññāyatana kusala citta'
$ws.Range('F74').Value = 'nevasaññānāsaññāyatana kusala citta'
$ws.Range('F75').Value = 'ākāsānañcāyatana vipāka citta'
$ws.Range('F76').Value = 'viññāṇañcāyatana vipāka citta'
$ws.Range('F77').Value = 'ākiñcaññāyatana vipāka citta'
$ws.Range('F78').Value = 'nevasaññānāsaññāyatana vipāka citta'
$ws.Range('F79').Value = 'ākāsānañcāyatana kiriya citta'
$ws.Range('F80').Value = 'viññāṇañcāyatana kiriya citta'
$ws.Range('F81').Value = 'ākiñcaññāyatana kiriya citta'
$ws.Range('F82').Value = 'nevasaññānāsaññāyatana kiriya citta'

# Widen new column to fit its content
$ws.Columns.Item(6).ColumnWidth = 52.3

# Make "Citta" sheet the active / selected sheet and tab
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range('J6').Select() | Out-Null

# Keep page orientation explicit (portrait)
$ws.PageSetup.Orientation = 1
